$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Style = "Normal"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.306.99'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +6.38%  '
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.990.22'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +3.68%  '
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.96%  '
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.79'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +6.30%  '
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("E7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("B8").Style = "Normal"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = 'XRP'
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Style = "Normal"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.515'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.13%  '
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("B9").Style = "Normal"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = 'LidoStakedEther'
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Style = "Normal"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.987.27'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +3.66%  '
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("E10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.63%  '
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("E11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +4.47%  '
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.447'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +3.51%  '
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000240'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.91%  '
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.85'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +5.78%  '
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("E15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.65%  '
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.304.47'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +6.39%  '
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.492.03'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +3.84%  '
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.90'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +5.23%  '
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.991.41'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.92%  '
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '450.57'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +4.43%  '
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.68'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +4.71%  '
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.680'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.66%  '
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.30'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +7.03%  '
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.16'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.32%  '
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.39'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +4.16%  '
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("B26").Style = "Normal"
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = 'RenderToken'
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Style = "Normal"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.75'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +7.34%  '
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("B27").Style = "Normal"
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'Fetch.AI'
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Style = "Normal"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.22'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +10.65%  '
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.09%  '
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.41'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +17.13%  '
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("E30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +11.48%  '
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("E31").Style = "Normal"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.70%  '
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("E32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.87%  '
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("E33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +4.02%  '
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.80'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +5.47%  '
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("E35").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.14%  '
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.983'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.43%  '
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.75'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +6.81%  '
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("E38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +9.16%  '
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '46.04'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +18.12%  '
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '49.12'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.51%  '
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.90'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.13%  '
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("E42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +5.83%  '
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.298'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +11.54%  '
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("E44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.37%  '
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '382.81'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +11.73%  '
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.765.21'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.07%  '
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("E47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +4.11%  '
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '134.60'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.12%  '
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("E49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.05%  '
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("E50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.59%  '
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.02'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +6.73%  '
$ws.Range("E51").Style = "Normal"
